$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (SFIA Level), shifting SFIA Level, Keycode,
# and Description one column to the right.
$ws.Columns("B:B").Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Skill Description"

# Map of SkillCode -> friendly Skill Description
$skillDescriptions = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "USEV"       = "User experience evaluation"
    "MADE"       = "MADE"
    "REQM"       = "Requirements definition and management"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -ne $null -and $code -ne "") {
        $ws.Cells.Item($r, 2).Value = $skillDescriptions[$code]
    }
}
